$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets: "Glaciers" -> "GIC", "AllButGlaciers" -> "All but GIC"
# ---------------------------------------------------------------------------
$wsGIC = $wb.Worksheets.Item("Glaciers")
$wsGIC.Name = "GIC"

$wsAllButGIC = $wb.Worksheets.Item("AllButGlaciers")
$wsAllButGIC.Name = "All but GIC"

# ---------------------------------------------------------------------------
# New shared strings need to be introduced in this order so they line up as
# uniqueCount index 89 ("#glac numbers from ...") then 90 ("1850-1900").
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2. "All but GIC" sheet (sheet9.xml): move the helper columns J:M -> Q:T,
#    fix up the formulas that reference them, and add the new GIC 1850-1900
#    row-10 data (S10, T10, V10).
# ---------------------------------------------------------------------------

# Row 4 header labels
$wsAllButGIC.Range("J4:M4").Copy()
$wsAllButGIC.Range("Q4").PasteSpecial(-4122)
$wsAllButGIC.Range("Q4").Value2 = $wsAllButGIC.Range("J4").Value2
$wsAllButGIC.Range("R4").Value2 = $wsAllButGIC.Range("K4").Value2
$wsAllButGIC.Range("S4").Value2 = $wsAllButGIC.Range("L4").Value2
$wsAllButGIC.Range("T4").Value2 = $wsAllButGIC.Range("M4").Value2
$wsAllButGIC.Range("J4:M4").Clear()

# Row 5
$wsAllButGIC.Range("J5:M5").Copy()
$wsAllButGIC.Range("Q5").PasteSpecial(-4122)
$wsAllButGIC.Range("Q5").Value2 = 1.35
$wsAllButGIC.Range("R5").Value2 = 0.34653615831230811
$wsAllButGIC.Range("S5").Value2 = 0.57999999999999996
$wsAllButGIC.Range("T5").Value2 = 0.14590996139465603
$wsAllButGIC.Range("J5:M5").Clear()
$wsAllButGIC.Range("E5").Formula = "=Q5-S5"
$wsAllButGIC.Range("F5").Formula = "=SQRT(R5^2+ T5^2)"

# Row 6
$wsAllButGIC.Range("J6:M6").Copy()
$wsAllButGIC.Range("Q6").PasteSpecial(-4122)
$wsAllButGIC.Range("Q6").Value2 = 2.33
$wsAllButGIC.Range("R6").Value2 = 0.47724716539502088
$wsAllButGIC.Range("S6").Value2 = 0.44
$wsAllButGIC.Range("T6").Value2 = 0.13983037966987874
$wsAllButGIC.Range("J6:M6").Clear()
$wsAllButGIC.Range("E6").Formula = "=Q6-S6"
$wsAllButGIC.Range("F6").Formula = "=SQRT(R6^2+ T6^2)"

# Row 7
$wsAllButGIC.Range("J7:M7").Copy()
$wsAllButGIC.Range("Q7").PasteSpecial(-4122)
$wsAllButGIC.Range("Q7").Value2 = 3.25
$wsAllButGIC.Range("R7").Value2 = 0.22190473295437274
$wsAllButGIC.Range("S7").Value2 = 0.55000000000000004
$wsAllButGIC.Range("T7").Value2 = 0.091193725871660011
$wsAllButGIC.Range("J7:M7").Clear()
$wsAllButGIC.Range("E7").Formula = "=Q7-S7"
$wsAllButGIC.Range("F7").Formula = "=SQRT(R7^2+ T7^2)"

# Row 8
$wsAllButGIC.Range("J8:M8").Copy()
$wsAllButGIC.Range("Q8").PasteSpecial(-4122)
$wsAllButGIC.Range("Q8").Value2 = 3.69
$wsAllButGIC.Range("R8").Value2 = 0.29181992278931213
$wsAllButGIC.Range("S8").Value2 = 0.62
$wsAllButGIC.Range("T8").Value2 = 0.033437699486275375
$wsAllButGIC.Range("J8:M8").Clear()
$wsAllButGIC.Range("E8").Formula = "=Q8-S8"
$wsAllButGIC.Range("F8").Formula = "=SQRT(R8^2+ T8^2)"

# Row 9
$wsAllButGIC.Range("J9:M9").Copy()
$wsAllButGIC.Range("Q9").PasteSpecial(-4122)
$wsAllButGIC.Range("Q9").Value2 = 1.73
$wsAllButGIC.Range("R9").Value2 = 0.27054138675259143
$wsAllButGIC.Range("S9").Value2 = 0.56999999999999995
$wsAllButGIC.Range("T9").Value2 = 0.13071100708271274
$wsAllButGIC.Range("J9:M9").Clear()
$wsAllButGIC.Range("E9").Formula = "=Q9-S9"
$wsAllButGIC.Range("F9").Formula = "=SQRT(R9^2+ T9^2)"

# Row 10 - J10/K10 keep their formula, move to Q10/R10; S10/T10/V10 are new
$wsAllButGIC.Range("J10:K10").Copy()
$wsAllButGIC.Range("Q10").PasteSpecial(-4122)
$wsAllButGIC.Range("Q10").Formula = "=0.014*1000/50"
$wsAllButGIC.Range("R10").Formula = "=0.014*1000/50"
$wsAllButGIC.Range("J10:K10").Clear()
$wsAllButGIC.Range("E10").Formula = "=Q10-S10"
$wsAllButGIC.Range("F10").Formula = "=SQRT(R10^2+ T10^2)"
$wsAllButGIC.Range("S10").Value2 = 0.434
$wsAllButGIC.Range("T10").Value2 = 0.11
$wsAllButGIC.Range("V10").Value2 = "#glac numbers from https://link.springer.com/article/10.1007/s10712-011-9121-7"

# ---------------------------------------------------------------------------
# 3. "GIC" sheet (sheet3.xml): add the new 1850-1900 row-10 estimate.
# ---------------------------------------------------------------------------
$wsGIC.Range("B9").Copy()
$wsGIC.Range("A10").PasteSpecial(-4122)
$wsGIC.Range("B9:C9").Copy()
$wsGIC.Range("B10").PasteSpecial(-4122)
$wsGIC.Range("A10").Value2 = "1850-1900"
$wsGIC.Range("B10").Value2 = 1850
$wsGIC.Range("C10").Value2 = 1900
$wsGIC.Range("E10").Value2 = 0.434
$wsGIC.Range("F10").Value2 = 0.11
$wsGIC.Range("G10").Value2 = "#glac numbers from https://link.springer.com/article/10.1007/s10712-011-9121-7"

# ---------------------------------------------------------------------------
# 4. GMSL sheet (sheet1.xml): A10 note changes from "PI" to "1850-1900".
# ---------------------------------------------------------------------------
$wsGMSL = $wb.Worksheets.Item("GMSL")
$wsGMSL.Range("A10").Value2 = "1850-1900"

# ---------------------------------------------------------------------------
# 5. Selections per sheet (also drives which tab is active / tabSelected).
# ---------------------------------------------------------------------------
$wsGMSL.Activate()
$wsGMSL.Range("G10").Select()

$wsTransposed = $wb.Worksheets.Item("transposed ar6")
$wsTransposed.Activate()
$wsTransposed.Range("C6").Select()

$wsAllButGIC.Activate()
$wsAllButGIC.Range("S10:W10").Select()

# GIC is the last-activated sheet, matching activeTab=2 / tabSelected on GIC.
$wsGIC.Activate()
$wsGIC.Range("D27").Select()
